# [GEN MCU] SCH 수정, PL 작성
#
# Applies:
#  - Sheet2!C5: Rup 10 -> 4.7 (ripples into C7 via existing formula)
#  - Sheet2: new LED current-calc block in columns H:J (rows 4-7)
#       H4 "LED"
#       H5 "V"   I5 3.3   J5 "V"
#       H6 "R"   I6 470   J6 "Ohm"
#       H7 "I"   I7 =I5/I6*1000   J7 "mA"
#  - Selection/cursor bookkeeping on Sheet1 & Sheet2 (cosmetic, matches
#    where the author last clicked before saving)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: SCH correction -> Rup changed from 10 to 4.7 ---
$ws2.Range("C5").Value = 4.7

# --- Sheet2: new PL (LED forward-current) calculation block ---
$ws2.Range("H4").Value = "LED"

$ws2.Range("H5").Value = "V"
$ws2.Range("I5").Value = 3.3
$ws2.Range("J5").Value = "V"

$ws2.Range("H6").Value = "R"
$ws2.Range("I6").Value = 470
$ws2.Range("J6").Value = "Ohm"

$ws2.Range("H7").Value = "I"
$ws2.Range("I7").Formula = "=I5/I6*1000"
$ws2.Range("I7").NumberFormat = $ws2.Range("D7").NumberFormat
$ws2.Range("J7").Value = "mA"

# --- restore the active-cell/selection bookkeeping Excel writes on save ---
[void]$ws1.Range("E25").Select()
[void]$ws2.Range("P23").Select()
[void]$ws2.Activate()
